# The deck's design theme ("Integral", carried in ppt/theme/theme1.xml and
# applied through the slide master) is replaced by the default Office
# color palette ("Office Theme"), matching the result of picking the
# built-in "Office Theme" design in the Design gallery. The font scheme
# and format scheme are already identical between the old and new theme,
# so only the twelve theme colors actually change.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# MsoThemeColorSchemeIndex order: 1=dk1 2=lt1 3=dk2 4=lt2 5-10=accent1-6
# 11=hlink 12=folHlink. Values are the Office Theme's RGB() (BGR-packed)
# equivalents of 000000 / FFFFFF / 44546A / E7E6E6 / 5B9BD5 / ED7D31 /
# A5A5A5 / FFC000 / 4472C4 / 70AD47 / 0563C1 / 954F72.
$colors.Item(1).RGB  = 0
$colors.Item(2).RGB  = 16777215
$colors.Item(3).RGB  = 6968388
$colors.Item(4).RGB  = 15132391
$colors.Item(5).RGB  = 13998939
$colors.Item(6).RGB  = 3243501
$colors.Item(7).RGB  = 10855845
$colors.Item(8).RGB  = 49407
$colors.Item(9).RGB  = 12874308
$colors.Item(10).RGB = 4697456
$colors.Item(11).RGB = 12673797
$colors.Item(12).RGB = 7491477
